$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): several "want-to-go" counts (column F) were incremented.
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F5").Value  = 126
$wsExh.Range("F10").Value = 16616
$wsExh.Range("F14").Value = 6441
$wsExh.Range("F21").Value = 62
$wsExh.Range("F31").Value = 5073
$wsExh.Range("F38").Value = 3854

# Sheet "全部类型" (sheet4): same events, mirrored, one row lower for the last item.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 126
$wsAll.Range("F10").Value = 16616
$wsAll.Range("F14").Value = 6441
$wsAll.Range("F21").Value = 62
$wsAll.Range("F31").Value = 5073
$wsAll.Range("F39").Value = 3854
